# Gap-Assessment.xlsx edit: rename sheet, refresh view scroll position, and
# re-apply the (slightly) recalculated row heights / column widths that Excel
# wrote out the next time the workbook was opened & saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab ---------------------------------------------
$ws.Name = "GAP"

# --- Column widths (character units) ---------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.0
$ws.Columns.Item(2).ColumnWidth = 18.833333333333332
$ws.Columns.Item(3).ColumnWidth = 10.5
$ws.Columns.Item(4).ColumnWidth = 19.833333333333332
$ws.Columns.Item(5).ColumnWidth = 18.0

# --- Row heights (points) ---------------------------------------------------
$ws.Rows.Item(1).RowHeight = 28.3
$ws.Rows.Item(2).RowHeight = 198
$ws.Rows.Item(3).RowHeight = 183.9
$ws.Rows.Item(4).RowHeight = 169.75
$ws.Rows.Item(5).RowHeight = 198
$ws.Rows.Item(6).RowHeight = 183.9
$ws.Rows.Item(7).RowHeight = 141.45
$ws.Rows.Item(8).RowHeight = 141.45
$ws.Rows.Item(9).RowHeight = 99
$ws.Rows.Item(10).RowHeight = 99
$ws.Rows.Item(11).RowHeight = 127.3
$ws.Rows.Item(12).RowHeight = 113.15
$ws.Rows.Item(13).RowHeight = 127.3
$ws.Rows.Item(14).RowHeight = 99
$ws.Rows.Item(15).RowHeight = 99
$ws.Rows.Item(16).RowHeight = 99
$ws.Rows.Item(17).RowHeight = 127.3
$ws.Rows.Item(18).RowHeight = 70.75
$ws.Rows.Item(19).RowHeight = 99
$ws.Rows.Item(20).RowHeight = 127.3
